$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 5000
$ws.Range("N4").Value = -5228
$ws.Range("H40").Value = 1792.2693
$ws.Range("I40").Value = 1659.96
$ws.Range("J40").Value = 5100
$ws.Range("K40").Value = 1659.96
$ws.Range("L40").Value = 5100
$ws.Range("M40").Value = -1484.96
$ws.Range("N40").Value = -5450
$ws.Range("H80").Value = 1500
$ws.Range("J80").Value = 1500
$ws.Range("L80").Value = 4500
$ws.Range("N80").Value = -6496
$ws.Range("H83").Value = 1500
$ws.Range("J83").Value = 1500
$ws.Range("L83").Value = 13500
$ws.Range("N83").Value = -23484
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H92").Value = 709.1429000000001
$ws.Range("I92").Value = 862.8
$ws.Range("K92").Value = 862.8
$ws.Range("M92").Value = 385.2
$ws.Range("H113").Value = 58333.332
$ws.Range("J113").Value = 58333.332
$ws.Range("L113").Value = 58333.332
$ws.Range("N113").Value = -64841.332
$ws.Range("H127").Value = 3333
$ws.Range("I127").Value = 2499.5
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 7498.5
$ws.Range("L127").Value = 15000
$ws.Range("M127").Value = -2538.5
$ws.Range("N127").Value = -24920
$ws.Range("H132").Value = 3858.4
$ws.Range("I132").Value = 1098
$ws.Range("K132").Value = 3294
$ws.Range("M132").Value = -764
$ws.Range("H138").Value = 2449.3333
$ws.Range("I138").Value = 939.2
$ws.Range("K138").Value = 2817.6
$ws.Range("M138").Value = 2322.4
$ws.Range("H141").Value = 3664
$ws.Range("I141").Value = 3478.5
$ws.Range("K141").Value = 10435.5
$ws.Range("M141").Value = -5255.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 709.5
$ws.Range("I2").Value = 696
$ws.Range("J2").Value = 750
$ws.Range("K2").Value = 696
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = -583
$ws.Range("N2").Value = -976
$ws.Range("H92").Value = 52400
$ws.Range("J92").Value = 52400
$ws.Range("L92").Value = 52400
$ws.Range("N92").Value = -57392
$ws.Range("H96").Value = 34668.8
$ws.Range("J96").Value = 36836
$ws.Range("L96").Value = 36836
$ws.Range("N96").Value = -42328
$ws.Range("H97").Value = 2856.5
$ws.Range("I97").Value = 2509.3333
$ws.Range("J97").Value = 3203.6667
$ws.Range("K97").Value = 2509.3333
$ws.Range("L97").Value = 3203.6667
$ws.Range("M97").Value = -2013.3333
$ws.Range("N97").Value = -4195.6667
$ws.Range("H116").Value = 709.5
$ws.Range("I116").Value = 696
$ws.Range("J116").Value = 750
$ws.Range("K116").Value = 696
$ws.Range("L116").Value = 750
$ws.Range("M116").Value = 1598
$ws.Range("N116").Value = -5338

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 709.5
$ws.Range("I3").Value = 696
$ws.Range("J3").Value = 750
$ws.Range("K3").Value = 696
$ws.Range("L3").Value = 750
$ws.Range("M3").Value = -582
$ws.Range("N3").Value = -978
$ws.Range("H94").Value = 2329.1667
$ws.Range("I94").Value = 1946.8182
$ws.Range("J94").Value = 2930
$ws.Range("K94").Value = 1946.8182
$ws.Range("L94").Value = 2930
$ws.Range("M94").Value = -1495.8182
$ws.Range("N94").Value = -3832
$ws.Range("H99").Value = 3758.3333
$ws.Range("I99").Value = 3537.5
$ws.Range("J99").Value = 4200
$ws.Range("K99").Value = 3537.5
$ws.Range("L99").Value = 4200
$ws.Range("M99").Value = -2039.5
$ws.Range("N99").Value = -7196

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 50864.25
$ws.Range("J43").Value = 50864.25
$ws.Range("L43").Value = 50864.25
$ws.Range("N43").Value = -51232.25
$ws.Range("H58").Value = 958.6923
$ws.Range("I58").Value = 534.625
$ws.Range("K58").Value = 534.625
$ws.Range("M58").Value = -331.625
$ws.Range("H99").Value = 1671166.6
$ws.Range("I99").Value = 2501000
$ws.Range("K99").Value = 2501000
$ws.Range("M99").Value = -2499502
$ws.Range("H101").Value = 50864.25
$ws.Range("J101").Value = 50864.25
$ws.Range("L101").Value = 50864.25
$ws.Range("N101").Value = -57354.25
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 1620
$ws.Range("I122").Value = 1212
$ws.Range("J122").Value = 1824
$ws.Range("K122").Value = 3636
$ws.Range("L122").Value = 5472
$ws.Range("M122").Value = -1186
$ws.Range("N122").Value = -10372
$ws.Range("H126").Value = 1671166.6
$ws.Range("I126").Value = 2501000
$ws.Range("K126").Value = 7503000
$ws.Range("M126").Value = -7500530
$ws.Range("H129").Value = 51997
$ws.Range("J129").Value = 51997
$ws.Range("L129").Value = 51997
$ws.Range("N129").Value = -61997
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 958.6923
$ws.Range("I136").Value = 534.625
$ws.Range("K136").Value = 1603.875
$ws.Range("M136").Value = 946.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1118.75
$ws.Range("I18").Value = 1307.3334
$ws.Range("J18").Value = 553
$ws.Range("K18").Value = 3922.0002
$ws.Range("L18").Value = 1659
$ws.Range("M18").Value = -3753.0002
$ws.Range("N18").Value = -1997
$ws.Range("H50").Value = 1168.4445
$ws.Range("I50").Value = 127.5
$ws.Range("J50").Value = 2001.2
$ws.Range("K50").Value = 382.5
$ws.Range("L50").Value = 6003.6
$ws.Range("M50").Value = 98.5
$ws.Range("N50").Value = -6965.6
$ws.Range("H53").Value = 1168.4445
$ws.Range("I53").Value = 127.5
$ws.Range("J53").Value = 2001.2
$ws.Range("K53").Value = 382.5
$ws.Range("L53").Value = 6003.6
$ws.Range("M53").Value = 98.5
$ws.Range("N53").Value = -6965.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11989.363
$ws.Range("I126").Value = 11188.3
$ws.Range("K126").Value = 33564.89999999999
$ws.Range("M126").Value = -31094.89999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2435.111
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 2902.6667
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 2902.6667
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -3492.6667
$ws.Range("H27").Value = 2435.111
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 2902.6667
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 2902.6667
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -3116.6667
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15580
$ws.Range("H46").Value = 252147.5
$ws.Range("I46").Value = 501499.75
$ws.Range("K46").Value = 501499.75
$ws.Range("M46").Value = -501311.75
$ws.Range("H136").Value = 3199.5
$ws.Range("I136").Value = 3199.5
$ws.Range("K136").Value = 9598.5
$ws.Range("M136").Value = -7048.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3200
$ws.Range("I126").Value = 2933.3333
$ws.Range("K126").Value = 8799.999899999999
$ws.Range("M126").Value = -6329.999899999999
$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200
